$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the two new daily rows (2026-01-13 / serial 46035) ---
# Copy the formatting from the last existing pair of rows (24:25 -> 四方坪站/高岭站)
# down onto the two new rows so the date / currency / integer number formats
# (and therefore the shared cellXfs entries) are reused instead of minting new styles.
$ws.Range("A24:F25").Copy()
$ws.Range("A26:F27").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Row 26: 四方坪站 (station name copied from an existing cell so the shared-string
# entry is reused rather than creating a new one)
$ws.Range("A26").Value2 = 46035
$ws.Range("B26").Value2 = $ws.Range("B2").Value2
$ws.Range("C26").Value2 = 13342.67
$ws.Range("D26").Value2 = 9743.2199999999993
$ws.Range("E26").Value2 = 3724.98
$ws.Range("F26").Value2 = 568

# Row 27: 高岭站
$ws.Range("A27").Value2 = 46035
$ws.Range("B27").Value2 = $ws.Range("B3").Value2
$ws.Range("C27").Value2 = 4654.72
$ws.Range("D27").Value2 = 3863.79
$ws.Range("E27").Value2 = 1337.27
$ws.Range("F27").Value2 = 172

# --- Update the saved view/selection to reflect scrolling down to the new rows ---
$ws.Activate() | Out-Null
$aw = $excel.ActiveWindow
$aw.ScrollRow = 13
$aw.ScrollColumn = 1
$ws.Range("L25").Select() | Out-Null
